$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header labels: OutSample -> TestSample
$ws.Range("E1").Value = "RMSE_TestSample"
$ws.Range("F1").Value = "R2_TestSample"
$ws.Range("G1").Value = "Adjusted_R2_TestSample"

# Update numeric values in E2:G9 with recomputed test-sample metrics
$ws.Range("E2").Value = 0.2203059453264854
$ws.Range("F2").Value = 0.8110024965460946
$ws.Range("G2").Value = 0.79350272770777

$ws.Range("E3").Value = 0.2227255026184086
$ws.Range("F3").Value = 0.7884419990383102
$ws.Range("G3").Value = 0.7644920366652888

$ws.Range("E4").Value = 0.219499513946479
$ws.Range("F4").Value = 0.8018701931147736
$ws.Range("G4").Value = 0.775198872957147

$ws.Range("E5").Value = 0.2056537017474372
$ws.Range("F5").Value = 0.8528421542171547
$ws.Range("G5").Value = 0.8297585705649437

$ws.Range("E6").Value = 0.207149563537544
$ws.Range("F6").Value = 0.8415319424109894
$ws.Range("G6").Value = 0.8130076920449675

$ws.Range("E7").Value = 0.2047829187007655
$ws.Range("F7").Value = 0.8430252898723016
$ws.Range("G7").Value = 0.8109896347441999

$ws.Range("E8").Value = 0.1908189957460074
$ws.Range("F8").Value = 0.8732255329845089
$ws.Range("G8").Value = 0.8441730509601255

$ws.Range("E9").Value = 0.1891215040397849
$ws.Range("F9").Value = 0.8780378382793181
$ws.Range("G9").Value = 0.8468985629463781
